$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (columns I and J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, center alignment) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for new columns I (I0) and J (IF), rows 2-15
$dataI = @(1,5,1,1,1,7,9,8,9,9,8,8,9,9)
$dataJ = @(1,6,1,2,2,8,9,8,9,9,9,8,9,9)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
